$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "V" marker for row 10 (matches the existing D1:D9 pattern)
$ws.Range("D10").Value = "V"

# Restore the active selection to F18 (single cell), as left by the author
$ws.Range("F18").Select()
